$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the completed data8 experiment row (row 11)
$ws.Range("C11").Value = 0.93879651623119498
$ws.Range("D11").Value = "tanh"
$ws.Range("E11").Value = 0.1
$ws.Range("F11").Value = "adam"
$ws.Range("G11").Value = 128

# Update the active selection to G11
$ws.Range("G11").Select()
